# Apply "hybrid bold + color" highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet/impact
# lines, matching the author's commit:
#   "Implement quantitative metrics highlighting across all resume formats"
#
# For DOCX the highlight is bold + dark slate color (#2C3E50). Each metric
# becomes its own run with <w:b/><w:color w:val="2C3E50"/>, while the
# surrounding text stays in plain runs - exactly what Word does when you
# select a sub-string of an existing run and toggle Bold/Font Color.

$BOLD_COLOR = 5258796   # RGB(0x2C, 0x3E, 0x50) -> w:color w:val="2C3E50"

# Bold+color every term in $terms, in order, within a single paragraph.
# Each search is confined to [cursor, paragraph-end) so it cannot leak
# into neighboring paragraphs and always advances past the previous hit,
# so repeated terms (e.g. "87%" then "71%") are matched left-to-right.
function Add-BoldHighlight($para, $terms) {
    $pEnd = $para.Range.End
    $cursor = $para.Range.Start
    foreach ($term in $terms) {
        $searchRng = $word.ActiveDocument.Range($cursor, $pEnd)
        $found = $searchRng.Find.Execute($term, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $searchRng.Font.Bold = 1
            $searchRng.Font.Color = $BOLD_COLOR
            $cursor = $searchRng.End
        }
    }
}

$d = $word.ActiveDocument

# "• Discovered systematic race coding errors ... from 23% to 64%"
Add-BoldHighlight $d.Paragraphs.Item(9) @("23%", "64%")

# "• Achieved 87% prediction accuracy ... from ±4.2% to ±2.1%"
Add-BoldHighlight $d.Paragraphs.Item(11) @("87%", "71%", "±4.2%", "±2.1%")

# "• Wrote RFP and analyzed bids from 1,200 vendors ..."
Add-BoldHighlight $d.Paragraphs.Item(31) @("1,200")

# "• Created comprehensive meta-analysis framework ... $400M ... $1B+"
Add-BoldHighlight $d.Paragraphs.Item(46) @("`$400M", "`$1B")

# "• Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Add-BoldHighlight $d.Paragraphs.Item(63) @("73.5%", "`$4.7M")

# "• Achieved 87% prediction accuracy ... industry standard of 71%"
Add-BoldHighlight $d.Paragraphs.Item(65) @("87%", "71%")
